$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin prices in column D are plain-looking numeric strings that must stay as
# text (matching the original inlineStr cells), so force text format before
# assigning the value - otherwise Excel auto-converts them to floating point
# numbers and e.g. "242.35" becomes 242.34999999999999 or drops trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.954.37"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.234.56"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.35"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.95"
$ws.Range("E7").Value = "  -0.58%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -4.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.24"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.91"
$ws.Range("E12").Value = "  -2.90%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.570.07"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.32"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.227.61"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.885.64"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("E19").Value = "  -6.65%  "

$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.27"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.05"
$ws.Range("E22").Value = "  +7.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.25"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("E24").Value = "  -5.94%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.32"
$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.45"
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.54"
$ws.Range("E31").Value = "  -1.89%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0796"
$ws.Range("E32").Value = "  -0.86%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.55"
$ws.Range("E33").Value = "  -4.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.68"
$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -6.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.29"
$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0303"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.64"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.44"
$ws.Range("E42").Value = "  +2.04%  "

$ws.Range("E43").Value = "  -1.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.68"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.76"
$ws.Range("E45").Value = "  -1.86%  "

$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("E49").Value = "  -2.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.69"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.445.77"
$ws.Range("E51").Value = "  +0.08%  "
